# Applies the coinranking.com snapshot refresh described in the commit
# "Updated cryptos list on Tue Jul  4 07:54:26 UTC 2023 with GitHub Actions":
# new Price (D) / Volume(1h) (E) readings for every coin row, plus the
# #48/#49 ranking swap between Algorand and EnergySwap (columns B-E).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.940.73"
$ws.Range("E2").Value = "  +0.97%  "
$ws.Range("D3").Value = "1.950.92"
$ws.Range("E3").Value = "  -0.71%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.18"
$ws.Range("E5").Value = "  -1.66%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4881"
$ws.Range("E7").Value = "  +1.68%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2961"
$ws.Range("E8").Value = "  +0.49%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06823"
$ws.Range("E9").Value = "  +0.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.13"
$ws.Range("E10").Value = "  -1.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "107.02"
$ws.Range("E11").Value = "  -4.57%  "
$ws.Range("D12").Value = "1.953.99"
$ws.Range("E12").Value = "  -0.51%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07743"
$ws.Range("E13").Value = "  +0.68%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.456"
$ws.Range("E14").Value = "  -0.69%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.7051"
$ws.Range("E15").Value = "  +2.37%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "281.02"
$ws.Range("D17").Value = "30.974.99"
$ws.Range("E17").Value = "  +1.02%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.22"
$ws.Range("E18").Value = "  -0.41%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007707"
$ws.Range("E19").Value = "  +0.39%  "
$ws.Range("D20").Value = "2.204.04"
$ws.Range("E20").Value = "  -0.72%  "
$ws.Range("E21").Value = "  +0.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.497"
$ws.Range("E22").Value = "  -2.99%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  +0.31%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.487"
$ws.Range("E24").Value = "  -1.85%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.836"
$ws.Range("E25").Value = "  +0.85%  "
$ws.Range("E26").Value = "  +0.26%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.93"
$ws.Range("E27").Value = "  -2.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.203"
$ws.Range("E28").Value = "  -0.35%  "
$ws.Range("E29").Value = "  -3.38%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.416"
$ws.Range("E30").Value = "  -1.26%  "
$ws.Range("E31").Value = "  -1.26%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.561"
$ws.Range("E32").Value = "  -2.70%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.460"
$ws.Range("E33").Value = "  +0.71%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04946"
$ws.Range("E34").Value = "  -2.71%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7695"
$ws.Range("E35").Value = "  -1.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.168"
$ws.Range("E36").Value = "  -0.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.730"
$ws.Range("E37").Value = "  +0.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02017"
$ws.Range("E38").Value = "  -2.65%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.700"
$ws.Range("E39").Value = "  -0.51%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.540"
$ws.Range("E40").Value = "  +8.40%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.154"
$ws.Range("E41").Value = "  +4.22%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "74.46"
$ws.Range("E42").Value = "  +6.82%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4489"
$ws.Range("E43").Value = "  +0.45%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "109.40"
$ws.Range("E44").Value = "  -1.45%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8831"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.117"
$ws.Range("E46").Value = "  +9.56%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9993"
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "978.24"
$ws.Range("E48").Value = "  +7.06%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.393"
$ws.Range("E49").Value = "  +0.69%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1266"
$ws.Range("E50").Value = "  +0.80%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "35.78"
$ws.Range("E51").Value = "  +0.15%  "
